$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two more Coverity-scan findings below the existing table (rows 6-7).
$ws.Range("A6").Value = "isThunderSecurityConfigured"
$ws.Range("B6").Value = "uninitialized variable"
$ws.Range("C6").Value = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/entservices-softwareupdate/helpers/UtilsController.h"
$ws.Range("D6").Value = 96

$ws.Range("A7").Value = "main"
$ws.Range("B7").Value = "uninitialized variable"
$ws.Range("C7").Value = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/devicesettings/generic/sample/getHostEDID.cpp"
$ws.Range("D7").Value = 49

$ws.Range("C14").Select()
